$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B14").Value = 1465
$ws.Range("C14").Value = 5300
$ws.Range("D14").Value = 21390
$ws.Range("E14").Value = 4919
$ws.Range("F14").Value = 6128
$ws.Range("G14").Value = 16229
$ws.Range("H14").Value = 79425
$ws.Range("I14").Value = 9694
$ws.Range("J14").Value = 7041
$ws.Range("K14").Value = 14650
$ws.Range("L14").Value = 5265
$ws.Range("M14").Value = 2575
$ws.Range("N14").Value = 6791
$ws.Range("O14").Value = 1068
$ws.Range("P14").Value = 1954
$ws.Range("Q14").Value = 183894
$ws.Range("R14").Value = 112
$ws.Range("S14").Value = 16507
